$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.693.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "'1.596.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'211.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.0619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'0.246"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'19.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "'1.820.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "'1.587.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "'65.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "'26.693.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'210.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -1.24%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'146.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'0.666"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.86%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "'1.294.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "'5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D44").Value = "'63.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "'1.732.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

# Rows 46 and 47: coin ranking order swapped (WEMIXToken now ranks above Quant)
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'0.882"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.56%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'90.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "

$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "'0.0999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "'0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "'7.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.40%  "
